$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.521.66'
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = '1.650.15'
$ws.Range("E3").Value = '  +1.06%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.001'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '300.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3783'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '50.88'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.46%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3566'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.16%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08112'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.80%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.224'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("E13").Value = '  -0.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.408'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.32%  '
$ws.Range("E15").Value = '  +0.87%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001202'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.54%  '
$ws.Range("D17").Value = '1.656.96'
$ws.Range("E17").Value = '  +2.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '97.16'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.34%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06990'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.805'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.31%  '
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("E23").Value = '  +1.47%  '
$ws.Range("D24").Value = '23.549.78'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.494'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.934'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.80%  '
$ws.Range("E27").Value = '  -0.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.99'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.242'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.52%  '
$ws.Range("D31").Value = '1.839.15'
$ws.Range("E31").Value = '  +2.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.962'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.78%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.158'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.45%  '
$ws.Range("E34").Value = '  +1.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.039'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02741'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("E38").Value = '  -1.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.976'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.47%  '
$ws.Range("E40").Value = '  +4.58%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.06894'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.99%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6931'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.320'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.63%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.67'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6445'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.13%  '
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("E47").Value = '  -1.28%  '
$ws.Range("E48").Value = '  -0.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07876'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '126.68'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.54%  '
$ws.Range("E51").Value = '  +0.20%  '
